# "Last fixes before exam"
# The only real content edit is the random-seed cell AE3 on the "Data"
# sheet. Everything else in the target diff (AK4, AF5:AF25) is a
# volatile INDIRECT()/MOD() formula chain that Excel recalculates
# automatically once the seed changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Activate()

# Update the seed value that drives the randomized lookup table.
$ws.Range("AE3").Value = 10566859

# Force a full recalculation so the dependent INDIRECT()/MOD() cells
# (AK4, AF5:AF25) pick up the new seed immediately.
$excel.CalculateFull()

# Reflect the author's final cursor position/selection on the sheet.
$ws.Range("AF8").Select() | Out-Null
